# Auto-generated edit script applying the diff to Sheets (Tonberry_Profits workbook)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 111 (sheet ALC)
$ws.Range("H111").Value = 1766.3334
$ws.Range("I111").Value = 1300
$ws.Range("J111").Value = 1999.5
$ws.Range("K111").Value = 3900
$ws.Range("L111").Value = 5998.5
$ws.Range("M111").Value = -833
$ws.Range("N111").Value = -12132.5

# Row 112 (sheet ALC)
$ws.Range("H112").Value = 3327.5715
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3327.5715
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 9982.7145
$ws.Range("N112").Value = -12198.7145

# Row 129 (sheet ALC)
$ws.Range("H129").Value = 972.09375
$ws.Range("I129").Value = 985
$ws.Range("J129").Value = 967.7917
$ws.Range("K129").Value = 2955
$ws.Range("L129").Value = 2903.3751
$ws.Range("M129").Value = 2045
$ws.Range("N129").Value = -12903.3751

# Row 138 (sheet ALC)
$ws.Range("H138").Value = 2763.7837
$ws.Range("I138").Value = 5178.1
$ws.Range("J138").Value = 1869.5927
$ws.Range("K138").Value = 15534.3
$ws.Range("L138").Value = 5608.7781
$ws.Range("M138").Value = -10394.3
$ws.Range("N138").Value = -15888.7781

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (sheet ARM)
$ws.Range("H61").Value = 3125.5454
$ws.Range("I61").Value = 1483
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 1483
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -1271
$ws.Range("N61").Value = -6424

# Row 122 (sheet ARM)
$ws.Range("H122").Value = 995.6667
$ws.Range("I122").Value = 995.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2987.0001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -537.0001000000002

# Row 136 (sheet ARM)
$ws.Range("H136").Value = 3125.5454
$ws.Range("I136").Value = 1483
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 4449
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -1899
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (sheet BSM)
$ws.Range("H99").Value = 1097.5
$ws.Range("I99").Value = 1097.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1097.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 400.5

# Row 105 (sheet BSM)
$ws.Range("H105").Value = 2545
$ws.Range("I105").Value = 2545
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2545
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -798
$ws.Range("N105").ClearContents()

# Row 108 (sheet BSM)
$ws.Range("H108").Value = 95000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 95000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 95000
$ws.Range("N108").Value = -102680

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (sheet CRP)
$ws.Range("H31").Value = 1475.2142
$ws.Range("I31").Value = 958.5714
$ws.Range("J31").Value = 1991.8572
$ws.Range("K31").Value = 958.5714
$ws.Range("L31").Value = 1991.8572
$ws.Range("M31").Value = -663.5714
$ws.Range("N31").Value = -2581.8572

# Row 34 (sheet CRP)
$ws.Range("H34").Value = 1475.2142
$ws.Range("I34").Value = 958.5714
$ws.Range("J34").Value = 1991.8572
$ws.Range("K34").Value = 958.5714
$ws.Range("L34").Value = 1991.8572
$ws.Range("M34").Value = -756.5714
$ws.Range("N34").Value = -2395.8572

# Row 70 (sheet CRP)
$ws.Range("H70").Value = 28833.334
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 28833.334
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 28833.334
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -29463.334

# Row 73 (sheet CRP)
$ws.Range("H73").Value = 28833.334
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 28833.334
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 28833.334
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -31017.334

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (sheet CUL)
$ws.Range("H2").Value = 555.3
$ws.Range("I2").Value = 550
$ws.Range("J2").Value = 556.625
$ws.Range("K2").Value = 3300
$ws.Range("L2").Value = 3339.75
$ws.Range("M2").Value = -3187
$ws.Range("N2").Value = -3565.75

# Row 7 (sheet CUL)
$ws.Range("H7").Value = 1521.2858
$ws.Range("I7").Value = 1521.2858
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4563.857400000001
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4451.857400000001
$ws.Range("N7").ClearContents()

# Row 23 (sheet CUL)
$ws.Range("H23").Value = 280.42856
$ws.Range("I23").Value = 295
$ws.Range("J23").Value = 278
$ws.Range("K23").Value = 885
$ws.Range("L23").Value = 834
$ws.Range("M23").Value = -650
$ws.Range("N23").Value = -1304

# Row 34 (sheet CUL)
$ws.Range("H34").Value = 329.23077
$ws.Range("I34").Value = 316.66666
$ws.Range("J34").Value = 480
$ws.Range("K34").Value = 949.9999799999999
$ws.Range("L34").Value = 1440
$ws.Range("M34").Value = -865.9999799999999
$ws.Range("N34").Value = -1608

# Row 68 (sheet CUL)
$ws.Range("H68").Value = 1942.1578
$ws.Range("I68").Value = 866.3333
$ws.Range("J68").Value = 2001.9259
$ws.Range("K68").Value = 2598.9999
$ws.Range("L68").Value = 6005.7777
$ws.Range("M68").Value = -1787.9999
$ws.Range("N68").Value = -7627.7777

# Row 71 (sheet CUL)
$ws.Range("H71").Value = 1942.1578
$ws.Range("I71").Value = 866.3333
$ws.Range("J71").Value = 2001.9259
$ws.Range("K71").Value = 7796.9997
$ws.Range("L71").Value = 18017.3331
$ws.Range("M71").Value = -3740.9997
$ws.Range("N71").Value = -26129.3331

# Row 116 (sheet CUL)
$ws.Range("H116").Value = 100002216
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 100002216
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 300006648
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -300013532

# Row 131 (sheet CUL)
$ws.Range("H131").Value = 19259034
$ws.Range("I131").Value = 50000404
$ws.Range("J131").Value = 45676.812
$ws.Range("K131").Value = 150001212
$ws.Range("L131").Value = 137030.436
$ws.Range("M131").Value = -149996172
$ws.Range("N131").Value = -147110.436

# Row 133 (sheet CUL)
$ws.Range("H133").Value = 125002500
$ws.Range("I133").Value = 250000000
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 750000000
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -749994940
$ws.Range("N133").Value = -25120

# Row 137 (sheet CUL)
$ws.Range("H137").Value = 4512.579
$ws.Range("I137").Value = 2363.3333
$ws.Range("J137").Value = 5504.5386
$ws.Range("K137").Value = 7089.999899999999
$ws.Range("L137").Value = 16513.6158
$ws.Range("M137").Value = -1989.999899999999
$ws.Range("N137").Value = -26713.6158

$ws = $wb.Worksheets.Item("GSM")
# Row 22 (sheet GSM)
$ws.Range("H22").Value = 41129.5
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 54339.332
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 54339.332
$ws.Range("M22").Value = -971
$ws.Range("N22").Value = -55397.332

# Row 122 (sheet GSM)
$ws.Range("H122").Value = 2564.7778
$ws.Range("I122").Value = 1451
$ws.Range("J122").Value = 3455.8
$ws.Range("K122").Value = 4353
$ws.Range("L122").Value = 10367.4
$ws.Range("M122").Value = -1903
$ws.Range("N122").Value = -15267.4

# Row 132 (sheet GSM)
$ws.Range("H132").Value = 7697707
$ws.Range("I132").Value = 19232268
$ws.Range("J132").Value = 7999.3335
$ws.Range("K132").Value = 57696804
$ws.Range("L132").Value = 23998.0005
$ws.Range("M132").Value = -57694274
$ws.Range("N132").Value = -29058.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (sheet LTW)
$ws.Range("H40").Value = 5945.4443
$ws.Range("I40").Value = 3358.4285
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 3358.4285
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -3222.4285
$ws.Range("N40").Value = -15272

# Row 133 (sheet LTW)
$ws.Range("H133").Value = 68738
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 68738
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 68738
$ws.Range("N133").Value = -73798

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (sheet WVR)
$ws.Range("H107").Value = 1701.25
$ws.Range("I107").Value = 1628.25
$ws.Range("J107").Value = 1847.25
$ws.Range("K107").Value = 4884.75
$ws.Range("L107").Value = 5541.75
$ws.Range("M107").Value = -2964.75
$ws.Range("N107").Value = -9381.75

# Row 126 (sheet WVR)
$ws.Range("H126").Value = 5612.3125
$ws.Range("I126").Value = 4780.7
$ws.Range("J126").Value = 6998.3335
$ws.Range("K126").Value = 14342.1
$ws.Range("L126").Value = 20995.0005
$ws.Range("M126").Value = -11872.1
$ws.Range("N126").Value = -25935.0005

# Row 128 (sheet WVR)
$ws.Range("H128").Value = 29078.572
$ws.Range("I128").Value = 10650
$ws.Range("J128").Value = 30000
$ws.Range("K128").Value = 10650
$ws.Range("L128").Value = 30000
$ws.Range("M128").Value = -5670
$ws.Range("N128").Value = -39960
